# feat(data): update demo data
#
# - episodios!E3 (Serie 1 / Episodio 2 "url"): point at the new coverr cat
#   video instead of the old bigbuckbunny sample.
# - showcase sheet: fill in the (until now empty) preview_img_url column
#   for both series, and replace the Serie 1 description placeholder with
#   real lorem-ipsum copy.
# - cosmetic: widen showcase!B a bit, and leave the "showcase" tab as the
#   active / selected one (it was "episodios" before).

$wb = $excel.ActiveWorkbook

$episodios = $wb.Worksheets.Item("episodios")
$showcase  = $wb.Worksheets.Item("showcase")

# --- episodios: new preview video url for Serie 1 / Episodio 2 ---------
$episodios.Range("E3").Value = "https://cdn.coverr.co/videos/coverr-a-beautiful-domestic-cat-with-green-eyes-4785/1080p.mp4"

# --- showcase: new copy + preview images --------------------------------
$showcase.Range("B2").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua."
$showcase.Range("C2").Value = "serie-1_showcase.jpg"
$showcase.Range("C3").Value = "serie-2_showcase.jpg"

# widen column B on showcase to fit the new description text
$showcase.Columns.Item(2).ColumnWidth = 16.1

# --- view state: episodios loses the selected tab, showcase gets it -----
[void]$episodios.Range("E16").Select()
[void]$showcase.Range("B2").Select()
[void]$showcase.Activate()
